# Update "想去人数" (F column) figures across the four worksheets to
# reflect the newly generated output (gh-pages regeneration at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 285
$ws1.Range("F4").Value  = 216
$ws1.Range("F5").Value  = 1306
$ws1.Range("F6").Value  = 211
$ws1.Range("F7").Value  = 595
$ws1.Range("F8").Value  = 120
$ws1.Range("F9").Value  = 579
$ws1.Range("F10").Value = 21
$ws1.Range("F11").Value = 643
$ws1.Range("F12").Value = 93
$ws1.Range("F14").Value = 152
$ws1.Range("F15").Value = 224

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 10
$ws2.Range("F6").Value = 2

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 6284
$ws3.Range("F4").Value = 1918

# --- Sheet 4: 全部类型 (All types, aggregate of the above) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 6284
$ws4.Range("F4").Value  = 1918
$ws4.Range("F5").Value  = 285
$ws4.Range("F9").Value  = 10
$ws4.Range("F11").Value = 216
$ws4.Range("F12").Value = 2
$ws4.Range("F15").Value = 1306
$ws4.Range("F17").Value = 211
$ws4.Range("F20").Value = 595
$ws4.Range("F22").Value = 120
$ws4.Range("F23").Value = 579
$ws4.Range("F24").Value = 21
$ws4.Range("F26").Value = 643
$ws4.Range("F27").Value = 93
$ws4.Range("F31").Value = 152
$ws4.Range("F37").Value = 224
